# Auto-generated: apply scraped-price updates to the Masamune_Profits workbook.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4793.615
$ws.Range("I15").Value = 4793.615
$ws.Range("K15").Value = 14380.845
$ws.Range("M15").Value = -14211.845
$ws.Range("H33").Value = 200.89394
$ws.Range("I33").Value = 159.11476
$ws.Range("K33").Value = 159.11476
$ws.Range("M33").Value = 69.88524000000001
$ws.Range("H64").Value = 27043.951
$ws.Range("I64").Value = 202240
$ws.Range("J64").Value = 2711.1667
$ws.Range("K64").Value = 202240
$ws.Range("L64").Value = 2711.1667
$ws.Range("M64").Value = -201992
$ws.Range("N64").Value = -3207.1667
$ws.Range("H67").Value = 27043.951
$ws.Range("I67").Value = 202240
$ws.Range("J67").Value = 2711.1667
$ws.Range("K67").Value = 202240
$ws.Range("L67").Value = 2711.1667
$ws.Range("M67").Value = -201382
$ws.Range("N67").Value = -4427.1667
$ws.Range("H113").Value = 2633.1516
$ws.Range("I113").Value = 2600.5715
$ws.Range("J113").Value = 2641.923
$ws.Range("K113").Value = 2600.5715
$ws.Range("L113").Value = 2641.923
$ws.Range("M113").Value = 653.4285
$ws.Range("N113").Value = -9149.922999999999
$ws.Range("H121").Value = 1181.3334
$ws.Range("J121").Value = 1494.091
$ws.Range("L121").Value = 4482.272999999999
$ws.Range("N121").Value = -7976.272999999999
$ws.Range("H132").Value = 25224.902
$ws.Range("I132").Value = 3952.6453
$ws.Range("J132").Value = 91168.89999999999
$ws.Range("K132").Value = 11857.9359
$ws.Range("L132").Value = 273506.7
$ws.Range("M132").Value = -9327.9359
$ws.Range("N132").Value = -278566.7
$ws.Range("H141").Value = 3247.04
$ws.Range("I141").Value = 1837.8
$ws.Range("J141").Value = 8884
$ws.Range("K141").Value = 5513.4
$ws.Range("L141").Value = 26652
$ws.Range("M141").Value = -333.3999999999996
$ws.Range("N141").Value = -37012

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32404.646
$ws.Range("I32").Value = 31804.646
$ws.Range("K32").Value = 31804.646
$ws.Range("M32").Value = -31517.646
$ws.Range("H45").Value = 1572.48
$ws.Range("I45").Value = 1397.5555
$ws.Range("K45").Value = 1397.5555
$ws.Range("M45").Value = -1020.5555

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2348.7856
$ws.Range("I107").Value = 2072.1667
$ws.Range("J107").Value = 2556.25
$ws.Range("K107").Value = 2072.1667
$ws.Range("L107").Value = 2556.25
$ws.Range("M107").Value = -152.1667000000002
$ws.Range("N107").Value = -6396.25
$ws.Range("H134").Value = 3919.8875
$ws.Range("I134").Value = 2977.4443
$ws.Range("K134").Value = 8932.332900000001
$ws.Range("M134").Value = -6397.332900000001

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1519.3889
$ws.Range("I22").Value = 289.92307
$ws.Range("K22").Value = 289.92307
$ws.Range("M22").Value = 60.07693
$ws.Range("H94").Value = 1939.6
$ws.Range("I94").Value = 1142.2
$ws.Range("J94").Value = 2737
$ws.Range("K94").Value = 1142.2
$ws.Range("L94").Value = 2737
$ws.Range("M94").Value = -691.2
$ws.Range("N94").Value = -3639
$ws.Range("H122").Value = 76084.5
$ws.Range("I122").Value = 150724.75
$ws.Range("J122").Value = 1444.25
$ws.Range("K122").Value = 452174.25
$ws.Range("L122").Value = 4332.75
$ws.Range("M122").Value = -449724.25
$ws.Range("N122").Value = -9232.75

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 166667470
$ws.Range("J80").Value = 333334500
$ws.Range("L80").Value = 1000003500
$ws.Range("N80").Value = -1000005372
$ws.Range("H83").Value = 166667470
$ws.Range("J83").Value = 333334500
$ws.Range("L83").Value = 3000010500
$ws.Range("N83").Value = -3000019860

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6832.28
$ws.Range("I70").Value = 7205.684
$ws.Range("J70").Value = 5649.8335
$ws.Range("K70").Value = 7205.684
$ws.Range("L70").Value = 5649.8335
$ws.Range("M70").Value = -6935.684
$ws.Range("N70").Value = -6189.8335
$ws.Range("H73").Value = 6832.28
$ws.Range("I73").Value = 7205.684
$ws.Range("J73").Value = 5649.8335
$ws.Range("K73").Value = 7205.684
$ws.Range("L73").Value = 5649.8335
$ws.Range("M73").Value = -6269.684
$ws.Range("N73").Value = -7521.8335
$ws.Range("H102").Value = 1416.1818
$ws.Range("I102").Value = 1530.8889
$ws.Range("K102").Value = 1530.8889
$ws.Range("M102").Value = 91.11110000000008

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2565.8064
$ws.Range("I7").Value = 2078.9092
$ws.Range("K7").Value = 2078.9092
$ws.Range("M7").Value = -1966.9092
$ws.Range("H61").Value = 3104.5454
$ws.Range("I61").Value = 2849.4443
$ws.Range("J61").Value = 4252.5
$ws.Range("K61").Value = 2849.4443
$ws.Range("L61").Value = 4252.5
$ws.Range("M61").Value = -2647.4443
$ws.Range("N61").Value = -4656.5
$ws.Range("H113").Value = 3104.5454
$ws.Range("I113").Value = 2849.4443
$ws.Range("J113").Value = 4252.5
$ws.Range("K113").Value = 2849.4443
$ws.Range("L113").Value = 4252.5
$ws.Range("M113").Value = -679.4443000000001
$ws.Range("N113").Value = -8592.5
$ws.Range("H126").Value = 2565.8064
$ws.Range("I126").Value = 2078.9092
$ws.Range("K126").Value = 6236.7276
$ws.Range("M126").Value = -3766.7276
$ws.Range("H136").Value = 2109.7
$ws.Range("I136").Value = 1649.75
$ws.Range("J136").Value = 3949.5
$ws.Range("K136").Value = 4949.25
$ws.Range("L136").Value = 11848.5
$ws.Range("M136").Value = -2399.25
$ws.Range("N136").Value = -16948.5

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3466.6667
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3466.6667
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3466.6667
$ws.Range("N62").Value = -4714.6667
$ws.Range("H65").Value = 3466.6667
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3466.6667
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 17333.3335
$ws.Range("N65").Value = -23573.3335
$ws.Range("H122").Value = 2041951.1
$ws.Range("J122").Value = 916.5
$ws.Range("L122").Value = 2749.5
$ws.Range("N122").Value = -7649.5
$ws.Range("H132").Value = 1827.8334
$ws.Range("I132").Value = 1337.875
$ws.Range("J132").Value = 2219.8
$ws.Range("K132").Value = 4013.625
$ws.Range("L132").Value = 6659.400000000001
$ws.Range("M132").Value = -1483.625
$ws.Range("N132").Value = -11719.4
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

